$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# New hole_id values for A2:A32, replacing the previous numeric index values
$holeIds = @(
    "BRG_16_04A",
    "BRG_13_01",
    "BRG_16_03",
    "BRG_05_11",
    "ECO_09_03",
    "BRG_01_06",
    "BRG_16_04B",
    "ECO_09_04",
    "BRG_01_02",
    "BRG_05_13",
    "BRG_01_03",
    "BRG_05_12",
    "BRG_05_09",
    "BRG_01_08",
    "BRG_05_04",
    "BRG_05_15",
    "ECO_09_02",
    "BRG_01_07",
    "BRG_13_02",
    "ECO_09_01",
    "BRG_16_08",
    "BRG_05_01",
    "BRG_16_02",
    "BRG_05_03",
    "BRG_05_02",
    "BRG_05_14",
    "BRG_08_01",
    "BRG_01_01",
    "BRG_01_09",
    "BRG_01_04",
    "BRG_16_01"
)

# Add the "hole_id" header in A1, copying the formatting already used by A2:A32
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "hole_id"

# Replace the numeric index values in A2:A32 with the new hole_id strings
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
